# Scotland league-two 2023-2024 — re-sync of match rows against the
# upstream scrape (rows reordered within a few matchdays + 5 new fixtures
# appended at the bottom; dimension grows from A1:V59 to A1:V64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row-rotation groups. Columns A:E (Indice/pais/torneio/temporada/
#    data_partida) stay put; F:V (the actual match record) gets moved
#    between rows within each group. Capture every source row's F:V
#    block BEFORE writing anything, so groups rotate safely.
# ---------------------------------------------------------------------

function Rotate-Rows {
    param($ws, [hashtable]$mapping)

    # Snapshot every row referenced (as source) first.
    $snapshots = @{}
    foreach ($srcRow in $mapping.Values) {
        if (-not $snapshots.ContainsKey($srcRow)) {
            $snapshots[$srcRow] = $ws.Range("F$($srcRow):V$($srcRow)").Value2
        }
    }

    foreach ($destRow in $mapping.Keys) {
        $srcRow = $mapping[$destRow]
        $ws.Range("F$($destRow):V$($destRow)").Value2 = $snapshots[$srcRow]
    }
}

# Rows 2-5 (05/08/2023 matchday)
Rotate-Rows $ws @{ 2 = 5; 3 = 2; 4 = 3; 5 = 4 }

# Rows 7-10 (12/08/2023 matchday)
Rotate-Rows $ws @{ 7 = 8; 8 = 9; 9 = 10; 10 = 7 }

# Rows 33-36 (23/09/2023 matchday)
Rotate-Rows $ws @{ 33 = 34; 34 = 35; 35 = 36; 36 = 33 }

# Rows 45-47 (21/10/2023 matchday)
Rotate-Rows $ws @{ 45 = 46; 46 = 47; 47 = 45 }

# ---------------------------------------------------------------------
# 2) Append 5 new fixtures (rows 60-64) for the 18/11/2023 matchday.
# ---------------------------------------------------------------------

$newRows = @(
    @{ A=59; E=45248.66666666666; F="Dumbarton";      G=1; H="Elgin City";      I=0;
       J=1.44; K="16/11/2023 09:13"; L=1.38; M="18/11/2023 15:52";
       N=4.29; O="16/11/2023 09:13"; P=4.86; Q="18/11/2023 15:52";
       R=5.68; S="16/11/2023 09:13"; T=7.89; U="18/11/2023 15:52";
       V="https://www.betexplorer.com/football/scotland/league-two/dumbarton-elgin-city/4ICTRinF/" },
    @{ A=60; E=45248.66666666666; F="Spartans";        G=1; H="Peterhead";       I=2;
       J=2.14; K="16/11/2023 09:13"; L=2.24; M="18/11/2023 15:41";
       N=3.35; O="16/11/2023 09:13"; P=3.57; Q="18/11/2023 15:41";
       R=2.98; S="16/11/2023 09:13"; T=3.02; U="18/11/2023 15:41";
       V="https://www.betexplorer.com/football/scotland/league-two/spartans-peterhead/xWAyQVGR/" },
    @{ A=61; E=45248.66666666666; F="Bonnyrigg Rose";  G=0; H="Forfar Athletic"; I=2;
       J=2.03; K="16/11/2023 09:13"; L=2.58; M="18/11/2023 15:58";
       N=3.33; O="16/11/2023 09:13"; P=3.09; Q="18/11/2023 15:52";
       R=3.23; S="16/11/2023 09:13"; T=2.88; U="18/11/2023 15:58";
       V="https://www.betexplorer.com/football/scotland/league-two/bonnyrigg-rose-forfar-athletic/pGdSTDH2/" },
    @{ A=62; E=45248.66666666666; F="Clyde";           G=0; H="East Fife";       I=4;
       J=2.85; K="16/11/2023 09:13"; L=2.86; M="18/11/2023 15:32";
       N=3.24; O="16/11/2023 09:13"; P=3.25; Q="18/11/2023 15:32";
       R=2.27; S="16/11/2023 09:13"; T=2.5;  U="18/11/2023 15:32";
       V="https://www.betexplorer.com/football/scotland/league-two/clyde-east-fife/O6eWSXW8/" },
    @{ A=63; E=45248.66666666666; F="Stranraer";       G=0; H="Stenhousemuir";   I=3;
       J=2.95; K="16/11/2023 09:13"; L=3.78; M="18/11/2023 15:18";
       N=3.33; O="16/11/2023 09:13"; P=3.49; Q="18/11/2023 15:58";
       R=2.17; S="16/11/2023 09:13"; T=1.97; U="18/11/2023 15:18";
       V="https://www.betexplorer.com/football/scotland/league-two/stranraer-stenhousemuir/r9BXQB1L/" }
)

$startRow = 60
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    # Pull formatting (bold/border/center for the index col, date format
    # for the match-date col) from the last existing data row so the
    # appended rows carry the same style indices as the rest of the sheet.
    $ws.Range("A59").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("E59").Copy() | Out-Null
    $ws.Range("E$row").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $data.A
    $ws.Cells.Item($row, 2).Value = "scotland"
    $ws.Cells.Item($row, 3).Value = "league-two"
    $ws.Cells.Item($row, 4).Value = "2023-2024"
    $ws.Cells.Item($row, 5).Value = $data.E
    $ws.Cells.Item($row, 6).Value = $data.F
    $ws.Cells.Item($row, 7).Value = $data.G
    $ws.Cells.Item($row, 8).Value = $data.H
    $ws.Cells.Item($row, 9).Value = $data.I
    $ws.Cells.Item($row, 10).Value = $data.J
    $ws.Cells.Item($row, 11).Value = $data.K
    $ws.Cells.Item($row, 12).Value = $data.L
    $ws.Cells.Item($row, 13).Value = $data.M
    $ws.Cells.Item($row, 14).Value = $data.N
    $ws.Cells.Item($row, 15).Value = $data.O
    $ws.Cells.Item($row, 16).Value = $data.P
    $ws.Cells.Item($row, 17).Value = $data.Q
    $ws.Cells.Item($row, 18).Value = $data.R
    $ws.Cells.Item($row, 19).Value = $data.S
    $ws.Cells.Item($row, 20).Value = $data.T
    $ws.Cells.Item($row, 21).Value = $data.U
    $ws.Cells.Item($row, 22).Value = $data.V
}
